$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New sequential "Ma SP" (product id) values for rows 2..23 (the 22 rows that survive).
$newIds = @("67","68","69","70","71","72","73","74","75","76","77","78","79","80","81","82","83","84","85","86","87","88")

# Force column A to store these as TEXT (matching the shared-string type used
# for the rest of the sheet) instead of Excel's default numeric auto-detection.
$idRange = $ws.Range("A2:A23")
$idRange.NumberFormat = "@"
for ($i = 0; $i -lt $newIds.Length; $i++) {
    $row = 2 + $i
    $ws.Range("A" + $row).Value = $newIds[$i]
}

# Restore the original cell style (border/font, no explicit number format) by
# copying the format from an untouched data cell back onto the range - this
# avoids leaving the range tagged with the temporary "@" text format.
$ws.Range("B2").Copy()
$idRange.PasteSpecial(-4122)

# Drop the last product ("Tra dao cam sa" / id 42) - row 24 - entirely.
$ws.Rows.Item(24).Delete()
